$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Delete the "Reffered_Physician" column (column H) from both sheets; the
# columns to its right (Co_Management, Anesthesiology_Physician, ...) shift
# left to fill the gap.
$ws1.Range("H1").EntireColumn.Delete()
$ws2.Range("H1").EntireColumn.Delete()

# Update the filter-database defined name so it no longer spans the deleted
# column (was $D$1:$M$3, now the sheet only goes to column L).
$wb.Names.Item("_xlnm._FilterDatabase").RefersTo = "='Nov 19 2020 - Dec 1 2020'!`$D`$1:`$L`$3"

# Sheet 1 is no longer the selected/active tab; sheet 2 becomes active and
# its selection moves to M2 (the new last column). Sheet 1's own selection
# also moves to M2.
$ws1.Range("M2").Select()

$ws2.Activate()
$ws2.Range("M2").Select()

Write-Host "done"
